$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = 2845956
$ws.Range("B14").Value = "Pril ISIS Cold Power liquid 3000ml Lemon"
$ws.Range("C13:D13").Copy()
$ws.Range("C14:D14").PasteSpecial(-4122)
$ws.Range("C13:D14").Select()
